$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 3.811642989160245),
    @(3, 3.182878228561681, 0.3375848360084654, 0.7127328510149897, 6.48142807727062, 0, 10.71462399285576),
    @(4, 0.02258322285507441, 0.3375848360084654, 3.082599426703578, 6.48142807727062, 0, 9.924195562837738),
    @(5, 1.505614041169197, 0.05231270169004087, 3.082599426703578, 6.48142807727062, 1, 11.12195424683344),
    @(6, 1.505614041169197, 9.226618575922256, 3.082599426703578, 6.48142807727062, 0, 20.29626012106565),
    @(7, 0.1554434735375247, 1766.335244827366, 0.7127328510149897, 246.9852506941017, 1, 2014.188671846021)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
